$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A27").Value = 26
$ws.Range("B27").Value = 26
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 27
$ws.Range("E27").Value = 36
$ws.Range("F27").Value = 54
$ws.Range("G27").Value = 90
